$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "retreat" area type: rows 1004/1005 were using the "UnpassArea"
# shared string by mistake for both the type (B) and name (C) columns;
# change them to the new "SpawnArea" value.
$ws.Range("B6").Value = "SpawnArea"
$ws.Range("C6").Value = "SpawnArea"
$ws.Range("B7").Value = "SpawnArea"
$ws.Range("C7").Value = "SpawnArea"

# Reset the sheet scroll position back to the top-left corner (A1) and
# move the active selection to H20.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
[void]$ws.Range("H20").Select()
